$d = $word.ActiveDocument

# Anchor on the last paragraph of the document body and make room for the
# new content by inserting a fresh paragraph mark right after it.
$lastPara = $d.Paragraphs.Last
$anchor = $lastPara.Range
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()

# The freshly created (still empty) paragraph is now the last paragraph;
# collapse its range to the start so the XML we inject lands *before* its
# paragraph mark instead of consuming it.
$newPara = $d.Paragraphs.Last
$insertionPoint = $newPara.Range
$insertionPoint.Collapse(1)

# Three new paragraphs: a blank spacer, an underlined "Common pitfalls:"
# heading, and a paragraph describing the pitfall. Supplying literal WordML
# guarantees the paragraphs do not inherit the numbered-list formatting of
# the preceding paragraph.
$openXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
      '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
          '<w:body>' + `
            '<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:cs="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p>' + `
            '<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:cs="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:cs="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Common pitfalls:</w:t></w:r></w:p>' + `
            '<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:cs="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:cs="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Confused about Shared and Distributed Memory concepts. </w:t></w:r><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p>' + `
          '</w:body>' + `
        '</w:document>' + `
      '</pkg:xmlData>' + `
    '</pkg:part>' + `
  '</pkg:package>'

$insertionPoint.InsertXML($openXml)

# Inserting the XML re-creates a trailing placeholder paragraph (a clone of
# the paragraph that used to be last, carrying its numbered-list pPr but no
# text). Remove it, together with the paragraph mark that precedes it, so
# the "Confused about..." paragraph becomes the document's true last
# paragraph, matching the target structure.
$trailing = $d.Paragraphs.Last
$cleanupRange = $d.Range($trailing.Range.Start - 1, $trailing.Range.End)
$cleanupRange.Delete()
